# Updated cryptos list on Sat Mar 25 02:54:05 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $val into $addr as plain TEXT (not auto-converted to a number)
# by building it via a formula in a scratch cell and pasting values-only,
# which avoids mutating number formats/styles on the target cell.
function Set-TextValue($addr, $val) {
    $ws.Range("Z1").Formula = '="' + $val + '"'
    $ws.Range("Z1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue "D2" "27.628.01"
$ws.Range("E2").Value = "  -2.34%  "

Set-TextValue "D3" "1.761.66"
$ws.Range("E3").Value = "  -3.14%  "

$ws.Range("E4").Value = "  -0.19%  "

Set-TextValue "D5" "324.72"
$ws.Range("E5").Value = "  -1.20%  "

Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  -0.08%  "

Set-TextValue "D7" "0.4275"
$ws.Range("E7").Value = "  -1.51%  "

Set-TextValue "D8" "0.3610"
$ws.Range("E8").Value = "  -2.03%  "

Set-TextValue "D9" "0.07582"
$ws.Range("E9").Value = "  -1.62%  "

Set-TextValue "D10" "42.37"
$ws.Range("E10").Value = "  -5.85%  "

Set-TextValue "D11" "1.107"
$ws.Range("E11").Value = "  -2.77%  "

Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  -0.08%  "

Set-TextValue "D13" "20.82"
$ws.Range("E13").Value = "  -6.16%  "

Set-TextValue "D14" "6.078"
$ws.Range("E14").Value = "  -3.96%  "

Set-TextValue "D15" "7.245"
$ws.Range("E15").Value = "  -4.06%  "

Set-TextValue "D16" "1.766.68"
$ws.Range("E16").Value = "  -3.98%  "

Set-TextValue "D17" "93.27"
$ws.Range("E17").Value = "  +0.19%  "

Set-TextValue "D18" "0.00001068"
$ws.Range("E18").Value = "  -1.43%  "

Set-TextValue "D19" "0.06435"
$ws.Range("E19").Value = "  -1.43%  "

Set-TextValue "D20" "0.9997"

Set-TextValue "D21" "17.18"
$ws.Range("E21").Value = "  -1.94%  "

Set-TextValue "D22" "5.899"
$ws.Range("E22").Value = "  -6.11%  "

Set-TextValue "D23" "27.674.32"
$ws.Range("E23").Value = "  -2.37%  "

Set-TextValue "D24" "11.31"
$ws.Range("E24").Value = "  -3.18%  "

Set-TextValue "D25" "2.125"
$ws.Range("E25").Value = "  +6.50%  "

Set-TextValue "D26" "162.50"
$ws.Range("E26").Value = "  +0.24%  "

Set-TextValue "D27" "20.40"
$ws.Range("E27").Value = "  -1.95%  "

Set-TextValue "D28" "1.963.32"
$ws.Range("E28").Value = "  -3.82%  "

Set-TextValue "D29" "2.158"
$ws.Range("E29").Value = "  -6.09%  "

Set-TextValue "D30" "125.52"
$ws.Range("E30").Value = "  -2.69%  "

Set-TextValue "D31" "1.107"
$ws.Range("E31").Value = "  -9.16%  "

Set-TextValue "D32" "5.613"
$ws.Range("E32").Value = "  -6.65%  "

Set-TextValue "D33" "3.666"
$ws.Range("E33").Value = "  +6.77%  "

Set-TextValue "D34" "0.08954"
$ws.Range("E34").Value = "  -2.75%  "

$ws.Range("E35").Value = "  -5.25%  "

Set-TextValue "D36" "0.02305"
$ws.Range("E36").Value = "  -2.39%  "

Set-TextValue "D37" "0.2114"
$ws.Range("E37").Value = "  -3.08%  "

Set-TextValue "D38" "0.06027"
$ws.Range("E38").Value = "  -2.78%  "

Set-TextValue "D39" "0.6365"
$ws.Range("E39").Value = "  -3.50%  "

Set-TextValue "D40" "4.959"
$ws.Range("E40").Value = "  -5.10%  "

$ws.Range("E41").Value = "  -0.49%  "

Set-TextValue "D42" "1.000"
$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E43").Value = "  -2.96%  "

Set-TextValue "D44" "7.912"
$ws.Range("E44").Value = "  -3.18%  "

Set-TextValue "D45" "13.49"
$ws.Range("E45").Value = "  -3.27%  "

Set-TextValue "D46" "0.5962"
$ws.Range("E46").Value = "  -2.68%  "

Set-TextValue "D47" "3.715"
$ws.Range("E47").Value = "  -1.15%  "

Set-TextValue "D48" "1.998"
$ws.Range("E48").Value = "  -1.40%  "

Set-TextValue "D49" "123.25"
$ws.Range("E49").Value = "  -2.33%  "

Set-TextValue "D50" "1.172"
$ws.Range("E50").Value = "  +1.10%  "

Set-TextValue "D51" "0.06867"
$ws.Range("E51").Value = "  -2.13%  "

$ws.Range("Z1").Clear()